$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 33 (ELE - GCYCP) entirely; rows below shift up
$ws.Rows.Item(33).Delete()

# Add new column K "Fecha" header
$ws.Range("A1").Copy()
$ws.Range("K1").PasteSpecial(-4122)
$ws.Range("K1").Value = "Fecha"

# Fill "Fecha" column for all data rows (2 through 42) with "07 08 24"
$ws.Range("K2:K42").Value = "07 08 24"
